$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '25.906.36'
$ws.Range('E2').Value = '  +0.18%  '
$ws.Range('D3').Value = '1.640.85'
$ws.Range('E3').Value = '  +0.49%  '
$ws.Range('D4').Value = '''1.007'
$ws.Range('E4').Value = '  -0.09%  '
$ws.Range('D5').Value = '''215.36'
$ws.Range('E5').Value = '  +0.11%  '
$ws.Range('D6').Value = '''0.5039'
$ws.Range('E6').Value = '  +0.34%  '
$ws.Range('D7').Value = '''1.006'
$ws.Range('E7').Value = '  -0.24%  '
$ws.Range('D8').Value = '''0.2575'
$ws.Range('E8').Value = '  +0.24%  '
$ws.Range('E9').Value = '  -0.59%  '
$ws.Range('D10').Value = '''19.55'
$ws.Range('E10').Value = '  +0.38%  '
$ws.Range('D11').Value = '''0.07805'
$ws.Range('E11').Value = '  +0.85%  '
$ws.Range('D12').Value = '1.663.04'
$ws.Range('E12').Value = '  +1.58%  '
$ws.Range('D13').Value = '''4.280'
$ws.Range('E13').Value = '  +0.78%  '
$ws.Range('D14').Value = '1.865.07'
$ws.Range('E14').Value = '  +0.41%  '
$ws.Range('E15').Value = '  -0.29%  '
$ws.Range('D16').Value = '0.0₅7886'
$ws.Range('E16').Value = '  -0.47%  '
$ws.Range('D17').Value = '''64.81'
$ws.Range('E17').Value = '  +2.19%  '
$ws.Range('D18').Value = '25.987.64'
$ws.Range('E18').Value = '  +0.53%  '
$ws.Range('E19').Value = '  -0.21%  '
$ws.Range('D20').Value = '''197.79'
$ws.Range('E20').Value = '  -2.74%  '
$ws.Range('D21').Value = '''4.389'
$ws.Range('E21').Value = '  +2.07%  '
$ws.Range('D22').Value = '''9.943'
$ws.Range('E22').Value = '  -0.36%  '
$ws.Range('D23').Value = '''5.980'
$ws.Range('E23').Value = '  +0.86%  '
$ws.Range('D24').Value = '''1.008'
$ws.Range('E24').Value = '  -0.16%  '
$ws.Range('D25').Value = '''1.881'
$ws.Range('E25').Value = '  -2.94%  '
$ws.Range('D26').Value = '''140.12'
$ws.Range('E26').Value = '  -0.67%  '
$ws.Range('D27').Value = '''0.1141'
$ws.Range('E27').Value = '  -0.97%  '
$ws.Range('D28').Value = '''6.850'
$ws.Range('E28').Value = '  +1.35%  '
$ws.Range('D29').Value = '''15.71'
$ws.Range('E29').Value = '  -0.20%  '
$ws.Range('D30').Value = '''1.242'
$ws.Range('E30').Value = '  +0.18%  '
$ws.Range('D31').Value = '''0.04868'
$ws.Range('E31').Value = '  -3.95%  '
$ws.Range('D32').Value = '''3.267'
$ws.Range('E32').Value = '  +0.28%  '
$ws.Range('D33').Value = '''3.196'
$ws.Range('E33').Value = '  +0.25%  '
$ws.Range('D34').Value = '''1.534'
$ws.Range('E34').Value = '  -0.63%  '
$ws.Range('E35').Value = '  +1.33%  '
$ws.Range('D36').Value = '''0.8900'
$ws.Range('E36').Value = '  -0.33%  '
$ws.Range('E37').Value = '  +0.40%  '
$ws.Range('B38').Value = 'Maker'
$ws.Range('C38').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D38').Value = '1.135.36'
$ws.Range('E38').Value = '  -0.24%  '
$ws.Range('B39').Value = 'ImmutableX'
$ws.Range('C39').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D39').Value = '''0.5541'
$ws.Range('E39').Value = '  -1.71%  '
$ws.Range('D40').Value = '''0.01564'
$ws.Range('E40').Value = '  +0.64%  '
$ws.Range('D41').Value = '''1.008'
$ws.Range('E41').Value = '  -0.11%  '
$ws.Range('D42').Value = '''5.699'
$ws.Range('E42').Value = '  +1.18%  '
$ws.Range('D43').Value = '''0.8174'
$ws.Range('E43').Value = '  +0.21%  '
$ws.Range('D44').Value = '''99.65'
$ws.Range('E44').Value = '  +0.35%  '
$ws.Range('D45').Value = '1.775.12'
$ws.Range('E45').Value = '  +0.37%  '
$ws.Range('E46').Value = '  +7.08%  '
$ws.Range('D47').Value = '''0.4533'
$ws.Range('E47').Value = '  +0.32%  '
$ws.Range('B48').Value = 'Frax'
$ws.Range('C48').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('D48').Value = '''1.010'
$ws.Range('E48').Value = '  -0.03%  '
$ws.Range('B49').Value = 'Aave'
$ws.Range('C49').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D49').Value = '''55.26'
$ws.Range('E49').Value = '  +0.89%  '
$ws.Range('E50').Value = '  +1.40%  '
$ws.Range('E51').Value = '  +0.01%  '
